# Update countries & provincias Spain
#
# The source "Pais" sheet is a COVID-19 ranking table sorted (descending)
# by "Casos totales" (column B). This edit refreshes the daily figures for
# a batch of countries. For two pairs of rows, the refreshed totals swap
# the row order of the two countries involved (Zambia now outranks
# Croacia; Islas Malvinas now outranks Montserrat), so those rows carry
# both a new country name (column A) and new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos -------------------------------------------------
$ws.Range("B4").Value = 6712072
$ws.Range("C4").Value = 3614
$ws.Range("D4").Value = 3975176
$ws.Range("E4").Value = 2538324
$ws.Range("G4").Value = 52
$ws.Range("H4").Value = 198572

# --- Row 19: Arabia Saudita -------------------------------------------------
$ws.Range("B19").Value = 326258
$ws.Range("C19").Value = 607
$ws.Range("D19").Value = 303930
$ws.Range("E19").Value = 18023
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 4305

# --- Row 25: Alemania --------------------------------------------------------
$ws.Range("B25").Value = 261895
$ws.Range("C25").Value = 597
$ws.Range("E25").Value = 16766
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 9429

# --- Row 49: Bielorrusia ------------------------------------------------------
$ws.Range("B49").Value = 74360
$ws.Range("C49").Value = 187
$ws.Range("D49").Value = 72609
$ws.Range("E49").Value = 995
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 756

# --- Row 60: Uzbekistan -------------------------------------------------------
$ws.Range("B60").Value = 47836
$ws.Range("C60").Value = 549
$ws.Range("D60").Value = 44375
$ws.Range("E60").Value = 3065
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 396

# --- Row 61: Suiza -------------------------------------------------------------
$ws.Range("E61").Value = 6511
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 2025

# --- Row 67: Azerbaiyan --------------------------------------------------------
$ws.Range("B67").Value = 38403
$ws.Range("C67").Value = 76
$ws.Range("D67").Value = 35860
$ws.Range("E67").Value = 1979
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 564

# --- Row 71: Serbia -------------------------------------------------------------
$ws.Range("B71").Value = 32437
$ws.Range("C71").Value = 29
$ws.Range("D71").Value = 31285
$ws.Range("E71").Value = 419

# --- Row 78: Bosnia y Herzegovina ------------------------------------------------
$ws.Range("B78").Value = 23635
$ws.Range("C78").Value = 170
$ws.Range("D78").Value = 16166
$ws.Range("E78").Value = 6764
$ws.Range("G78").Value = 9
$ws.Range("H78").Value = 705

# --- Row 79: Libia ----------------------------------------------------------------
$ws.Range("B79").Value = 23515
$ws.Range("C79").Value = 734
$ws.Range("D79").Value = 12762
$ws.Range("E79").Value = 10385
$ws.Range("G79").Value = 6
$ws.Range("H79").Value = 368

# --- Row 88/89: Zambia overtakes Croacia -------------------------------------------
$ws.Range("A88").Value = "Zambia"
$ws.Range("B88").Value = 13720
$ws.Range("C88").Value = 181
$ws.Range("D88").Value = 12380
$ws.Range("E88").Value = 1020
$ws.Range("G88").Value = 8
$ws.Range("H88").Value = 320

$ws.Range("A89").Value = "Croacia"
$ws.Range("B89").Value = 13598
$ws.Range("C89").Value = 65
$ws.Range("D89").Value = 11151
$ws.Range("E89").Value = 2220
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 227

# --- Row 98: Namibia ------------------------------------------------------------
$ws.Range("B98").Value = 9818
$ws.Range("C98").Value = 99
$ws.Range("D98").Value = 6693
$ws.Range("E98").Value = 3022
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 103

# --- Row 104: Haiti ---------------------------------------------------------------
$ws.Range("B104").Value = 8499
$ws.Range("C104").Value = 6
$ws.Range("E104").Value = 2160

# --- Row 150: Islandia -------------------------------------------------------------
$ws.Range("B150").Value = 2168
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 2095

# --- Row 162: Liberia --------------------------------------------------------------
$ws.Range("B162").Value = 1321
$ws.Range("C162").Value = 2
$ws.Range("D162").Value = 1213
$ws.Range("E162").Value = 26

# --- Row 214/215: Islas Malvinas overtakes Montserrat -------------------------------
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Timestamp footer ---------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 16:23"
